$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.42612573704929
$ws.Range("C2").Value = 0.2716980152002861
$ws.Range("D2").Value = 0.02815277586734055
$ws.Range("E2").Value = 0.09735005301324939
$ws.Range("F2").Value = 0.7404525498734955
$ws.Range("L2").Value = 0.2186133324869814
$ws.Range("O2").Value = 2.54532475422431
$ws.Range("B3").Value = 1.287646044421081
$ws.Range("C3").Value = 0.2566188330006298
$ws.Range("D3").Value = 0.02670168647996718
$ws.Range("E3").Value = 0.09856702563687314
$ws.Range("F3").Value = 0.7407862203639866
$ws.Range("L3").Value = 0.2082890078714001
$ws.Range("O3").Value = 2.562567931480373
$ws.Range("B4").Value = 1.202640181682398
$ws.Range("C4").Value = 0.2473352152663608
$ws.Range("D4").Value = 0.02580839362537546
$ws.Range("E4").Value = 0.09937666743665829
$ws.Range("F4").Value = 0.7416464198107704
$ws.Range("L4").Value = 0.2020494733446441
$ws.Range("O4").Value = 2.575409500003616
$ws.Range("B5").Value = 1.168006751467431
$ws.Range("C5").Value = 0.2435460889222441
$ws.Range("D5").Value = 0.02544381431743403
$ws.Range("E5").Value = 0.09972228918227444
$ws.Range("F5").Value = 0.7421613517768151
$ws.Range("L5").Value = 0.1995319262777286
$ws.Range("O5").Value = 2.581207977441608
$ws.Range("B6").Value = 1.162256382369662
$ws.Range("C6").Value = 0.2429165551417043
$ws.Range("D6").Value = 0.02538324349059451
$ws.Range("E6").Value = 0.09978062655592979
$ws.Range("F6").Value = 0.7422567732847014
$ws.Range("L6").Value = 0.1991154080681099
$ws.Range("O6").Value = 2.582204919958031
$ws.Range("B7").Value = 1.202173071766197
$ws.Range("C7").Value = 0.2472841375883092
$ws.Range("D7").Value = 0.02580347899241531
$ws.Range("E7").Value = 0.09938126511149292
$ws.Range("F7").Value = 0.7416526992634473
$ws.Range("L7").Value = 0.202015419060416
$ws.Range("O7").Value = 2.57548541274511
$ws.Range("B8").Value = 1.378374569436744
$ws.Range("C8").Value = 0.2665040599692077
$ws.Range("D8").Value = 0.02765293741340002
$ws.Range("E8").Value = 0.09775670556897786
$ws.Range("F8").Value = 0.7404313483539582
$ws.Range("L8").Value = 0.2150328553357923
$ws.Range("O8").Value = 2.550801499978348
$ws.Range("B9").Value = 1.724012180241402
$ws.Range("C9").Value = 0.3039854258961441
$ws.Range("D9").Value = 0.03126030033720184
$ws.Range("E9").Value = 0.09506657982688083
$ws.Range("F9").Value = 0.7432548820246581
$ws.Range("L9").Value = 0.2413496595601003
$ws.Range("O9").Value = 2.520347592924026
$ws.Range("B10").Value = 1.977961545731205
$ws.Range("C10").Value = 0.3313838256891302
$ws.Range("D10").Value = 0.03389769180461144
$ws.Range("E10").Value = 0.09339282582802255
$ws.Range("F10").Value = 0.748538571820859
$ws.Range("L10").Value = 0.2611670464070528
$ws.Range("O10").Value = 2.509009849006873
$ws.Range("B11").Value = 2.09348190978568
$ws.Range("C11").Value = 0.3438155814336028
$ws.Range("D11").Value = 0.03509449064628711
$ws.Range("E11").Value = 0.09269723182207201
$ws.Range("F11").Value = 0.7516451862771305
$ws.Range("L11").Value = 0.2702876302099781
$ws.Range("O11").Value = 2.506268740626155
$ws.Range("B12").Value = 2.137224739684314
$ws.Range("C12").Value = 0.348518330973377
$ws.Range("D12").Value = 0.03554723832883866
$ws.Range("E12").Value = 0.09244330170099069
$ws.Range("F12").Value = 0.7529231317853515
$ws.Range("L12").Value = 0.2737565190597593
$ws.Range("O12").Value = 2.5055798245331
$ws.Range("B13").Value = 2.127804067047634
$ws.Range("C13").Value = 0.3475057317340315
$ws.Range("D13").Value = 0.03544975164790998
$ws.Range("E13").Value = 0.09249756842254797
$ws.Range("F13").Value = 0.7526433804024322
$ws.Range("L13").Value = 0.273008760134303
$ws.Range("O13").Value = 2.505712644077278
$ws.Range("B14").Value = 2.097080713932939
$ws.Range("C14").Value = 0.3442025790894832
$ws.Range("D14").Value = 0.0351317477288049
$ws.Range("E14").Value = 0.09267615084272052
$ws.Range("F14").Value = 0.7517482861240268
$ws.Range("L14").Value = 0.2705727150445796
$ws.Range("O14").Value = 2.506205057842948
$ws.Range("B15").Value = 2.078261428251494
$ws.Range("C15").Value = 0.3421786575941042
$ws.Range("D15").Value = 0.03493690119277915
$ws.Range("E15").Value = 0.09278677221841392
$ws.Range("F15").Value = 0.7512132513722634
$ws.Range("L15").Value = 0.2690825347070245
$ws.Range("O15").Value = 2.506552183530147
$ws.Range("B16").Value = 1.970411808117717
$ws.Range("C16").Value = 0.3305707135021123
$ws.Range("D16").Value = 0.03381941615430861
$ws.Range("E16").Value = 0.0934396098567909
$ws.Range("F16").Value = 0.74834973200565
$ws.Range("L16").Value = 0.2605731139490075
$ws.Range("O16").Value = 2.50923774917436
$ws.Range("B17").Value = 1.904247587643169
$ws.Range("C17").Value = 0.3234412288701378
$ws.Range("D17").Value = 0.03313309620898508
$ws.Range("E17").Value = 0.09385696850501368
$ws.Range("F17").Value = 0.7467734355434317
$ws.Range("L17").Value = 0.255379854243742
$ws.Range("O17").Value = 2.511505312961702
$ws.Range("B18").Value = 1.866191555075602
$ws.Range("C18").Value = 0.3193375514979948
$ws.Range("D18").Value = 0.03273806580278915
$ws.Range("E18").Value = 0.09410321523863097
$ws.Range("F18").Value = 0.7459329384070799
$ws.Range("L18").Value = 0.2524027673108407
$ws.Range("O18").Value = 2.513036911941725
$ws.Range("B19").Value = 1.853306472577458
$ws.Range("C19").Value = 0.3179476130713681
$ws.Range("D19").Value = 0.03260426851982601
$ws.Range("E19").Value = 0.0941876534057382
$ws.Range("F19").Value = 0.7456597078140561
$ws.Range("L19").Value = 0.2513964855611164
$ws.Range("O19").Value = 2.5135944850675
$ws.Range("B20").Value = 1.911290908205331
$ws.Range("C20").Value = 0.3242004859288556
$ws.Range("D20").Value = 0.03320618504464079
$ws.Range("E20").Value = 0.09381189890589958
$ws.Range("F20").Value = 0.7469343858267479
$ws.Range("L20").Value = 0.2559316574051422
$ws.Range("O20").Value = 2.511240384834366
$ws.Range("B21").Value = 2.10610498263992
$ws.Range("C21").Value = 0.3451729300581121
$ws.Range("D21").Value = 0.03522516573176659
$ws.Range("E21").Value = 0.09262343956701713
$ws.Range("F21").Value = 0.7520084377834593
$ws.Range("L21").Value = 0.2712878305882498
$ws.Range("O21").Value = 2.506050937310818
$ws.Range("B22").Value = 2.233413424288301
$ws.Range("C22").Value = 0.3588510172441488
$ws.Range("D22").Value = 0.03654202706866272
$ws.Range("E22").Value = 0.09190195369455267
$ws.Range("F22").Value = 0.7559166693299488
$ws.Range("L22").Value = 0.281412092657348
$ws.Range("O22").Value = 2.504694788128177
$ws.Range("B23").Value = 2.165468388404463
$ws.Range("C23").Value = 0.3515534802882314
$ws.Range("D23").Value = 0.03583944580751108
$ws.Range("E23").Value = 0.09228196525756438
$ws.Range("F23").Value = 0.7537764548151955
$ws.Range("L23").Value = 0.2760005394590337
$ws.Range("O23").Value = 2.50523182063381
$ws.Range("B24").Value = 1.908106676857642
$ws.Range("C24").Value = 0.3238572408721723
$ws.Range("D24").Value = 0.0331731429993809
$ws.Range("E24").Value = 0.09383225524341832
$ws.Range("F24").Value = 0.7468614154691267
$ws.Range("L24").Value = 0.2556821604460282
$ws.Range("O24").Value = 2.511359448910014
$ws.Range("B25").Value = 1.630502460954688
$ws.Range("C25").Value = 0.293869323278642
$ws.Range("D25").Value = 0.03028660745813028
$ws.Range("E25").Value = 0.09574121442095063
$ws.Range("F25").Value = 0.7419294698460135
$ws.Range("L25").Value = 0.2341456260844836
$ws.Range("O25").Value = 2.526654732982166
